$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: LinearRegression - slight precision changes
$ws.Range("B2").Value = 0.6644204749112232
$ws.Range("C2").Value = 0.6644204749112232
$ws.Range("D2").Value = 0.6644204749112232

# Row 3: RandomForestRegressor - value changes
$ws.Range("B3").Value = 0.9155408654331411
$ws.Range("C3").Value = 0.9297152915404864
$ws.Range("D3").Value = 0.9237118633856677

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, value changes
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9561894533603632
$ws.Range("C4").Value = 0.9226174920391154
$ws.Range("D4").Value = 0.956953168454318

# Row 5: AdaBoostRegressor -> MLPRegressor, value changes
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.5130151725227473
$ws.Range("C5").Value = 0.5136278684860438
$ws.Range("D5").Value = 0.4798668217832693
